# Generate Report for Handoff
# The cbf32fc4-... source file (row 3 on the Overview / zh-cn / de-de sheets)
# has now been handed off for localization. Update its status from
# "Handed back: in sync with en-US" to "Ready for handoff" and stamp the
# new "Latest Handoff Datetime" for each language.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 reflects the per-language status for the
#     cbf32fc4 file in both the zh-cn (col B) and de-de (col C) columns.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row 3 is the cbf32fc4 file.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-03-08 06:24:29"

# --- de-de sheet: row 3 is the cbf32fc4 file.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-03-08 06:24:32"
